$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.352.80'
$ws.Range("E2").Value = '  +7.21%  '
$ws.Range("D3").Value = '3.647.21'
$ws.Range("E3").Value = '  +4.51%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '421.01'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.81'
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.667'
$ws.Range("E7").Value = '  +5.72%  '
$ws.Range("D8").Value = '3.639.66'
$ws.Range("E8").Value = '  +4.48%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.778'
$ws.Range("E10").Value = '  +5.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.197'
$ws.Range("E11").Value = '  +29.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000420'
$ws.Range("E12").Value = '  +89.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.47'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.95'
$ws.Range("E14").Value = '  +1.75%  '
$ws.Range("D15").Value = '4.211.10'
$ws.Range("E15").Value = '  +4.30%  '
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.30'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '3.629.28'
$ws.Range("E18").Value = '  +5.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.14'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '68.301.82'
$ws.Range("E20").Value = '  +7.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.56'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '460.73'
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.82'
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.10'
$ws.Range("E24").Value = '  -5.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.43'
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.22'
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.15'
$ws.Range("E28").Value = '  +7.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.87'
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.79'
$ws.Range("E30").Value = '  +4.53%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.46'
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  +3.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.26'
$ws.Range("E33").Value = '  -3.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.158'
$ws.Range("E34").Value = '  -5.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.33'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0839'
$ws.Range("E36").Value = '  +34.58%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.19'
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0489'
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.151'
$ws.Range("E40").Value = '  +10.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '148.65'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.97'
$ws.Range("E43").Value = '  -5.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.69'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.30'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.30'
$ws.Range("E46").Value = '  -7.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.171'
$ws.Range("E47").Value = '  +21.69%  '
$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.58'
$ws.Range("E48").Value = '  +8.75%  '
$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.305'
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.97'
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.64'
$ws.Range("E51").Value = '  +14.33%  '
